{"js": "// The document body that the script will operate on.\nconst body = context.document.body;\n\n// ------------------------------------------------------------------\n// 1) \"Modify\" -> \"Modify: \"\n//    (the list item that used to just say \"Modify\" now also has a\n//    trailing \": \" \u2014 reflects the commit's note that Modify now takes\n//    a request body, described right after the colon).\n// ------------------------------------------------------------------\nconst modifyResults = body.search(\"Modify\", { matchCase: true, matchWholeWord: true });\nmodifyResults.load(\"items\");\nawait context.sync();\n\nif (modifyResults.items.length > 0) {\n  // Insert right after the found \"Modify\" text.\n  modifyResults.items[0].insertText(\": \", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// ------------------------------------------------------------------\n// 2) Move the \"_GoBack\" bookmark from the very end of the document to\n//    the paragraph holding the two screenshots that illustrate the\n//    \"Modify\" endpoint (right before the first of those two images) \u2014\n//    this mirrors where Word leaves \"_GoBack\" after the last edit was\n//    made in that paragraph.\n// ------------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the \"Modify\" list paragraph, then the very next paragraph is the\n// one that contains its two screenshots.\nlet modifyParagraphIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Modify\") === 0) {\n    modifyParagraphIndex = i;\n    break;\n  }\n}\n\nif (modifyParagraphIndex !== -1 && modifyParagraphIndex + 1 < paragraphs.items.length) {\n  const imagesParagraph = paragraphs.items[modifyParagraphIndex + 1];\n  const startOfParagraph = imagesParagraph.getRange(Word.RangeLocation.start);\n  startOfParagraph.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) \"Modify\" -> \"Modify: \"\n#    (the list item that used to just say \"Modify\" now also has a\n#    trailing \": \" \u2014 reflects the commit's note that Modify now takes\n#    a request body, described right after the colon).\n# ------------------------------------------------------------------\n$find = $d.Content.Find\n$find.Text = \"Modify\"\n$find.MatchWholeWord = $true\n$find.Execute(\"Modify\", $false, $true, $false, $false, $false, $true, 1, $false, \"Modify: \", 2) | Out-Null\n\n# ------------------------------------------------------------------\n# 2) Move the \"_GoBack\" bookmark from the very end of the document to\n#    the paragraph holding the two screenshots that illustrate the\n#    \"Modify\" endpoint (right before the first of those two images) \u2014\n#    this mirrors where Word leaves \"_GoBack\" after the last edit was\n#    made in that paragraph.\n# ------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Find the \"Modify\" list paragraph, then the very next paragraph is the\n# one that contains its two screenshots.\n$modifyParaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"Modify\")) {\n        $modifyParaIndex = $i\n        break\n    }\n}\n\nif ($modifyParaIndex -ne -1 -and ($modifyParaIndex + 1) -le $d.Paragraphs.Count) {\n    $imagesPara = $d.Paragraphs.Item($modifyParaIndex + 1)\n    $d.Bookmarks.Add(\"_GoBack\", $imagesPara.Range) | Out-Null\n}\n"}
